# Updates cryptos list values (Price / Volume(1h) columns, and a row
# re-ordering of BabyDogeCoin/RenderToken) per the "Updated cryptos list"
# GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.912.74'
$ws.Range("E2").Value = '  -4.66%  '

# Row 3
$ws.Range("D3").Value = '1.737.69'
$ws.Range("E3").Value = '  -5.07%  '

# Row 4
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.29%  '

# Row 5
$ws.Range("D5").Value = '''226.35'
$ws.Range("E5").Value = '  -3.97%  '

# Row 6
$ws.Range("D6").Value = '''0.5784'
$ws.Range("E6").Value = '  -3.90%  '

# Row 7
$ws.Range("E7").Value = '  -0.29%  '

# Row 8
$ws.Range("D8").Value = '''0.2728'
$ws.Range("E8").Value = '  -1.57%  '

# Row 9
$ws.Range("D9").Value = '''23.23'
$ws.Range("E9").Value = '  -1.37%  '

# Row 10
$ws.Range("D10").Value = '''0.06608'
$ws.Range("E10").Value = '  -5.21%  '

# Row 11
$ws.Range("D11").Value = '''0.07543'
$ws.Range("E11").Value = '  -1.01%  '

# Row 12
$ws.Range("D12").Value = '1.741.13'
$ws.Range("E12").Value = '  -5.80%  '

# Row 13
$ws.Range("E13").Value = '  -1.09%  '

# Row 14
$ws.Range("D14").Value = '''0.6020'
$ws.Range("E14").Value = '  -4.75%  '

# Row 15
$ws.Range("D15").Value = '1.973.51'
$ws.Range("E15").Value = '  -5.11%  '

# Row 16
$ws.Range("D16").Value = '''74.58'
$ws.Range("E16").Value = '  -4.11%  '

# Row 17
$ws.Range("D17").Value = '''0.000008752'
$ws.Range("E17").Value = '  -11.38%  '

# Row 18
$ws.Range("D18").Value = '27.915.08'
$ws.Range("E18").Value = '  -3.61%  '

# Row 19
$ws.Range("D19").Value = '''5.322'
$ws.Range("E19").Value = '  -4.72%  '

# Row 20
$ws.Range("E20").Value = '  -0.28%  '

# Row 21
$ws.Range("D21").Value = '''205.43'
$ws.Range("E21").Value = '  -5.53%  '

# Row 22
$ws.Range("E22").Value = '  -2.78%  '

# Row 23
$ws.Range("D23").Value = '''6.613'
$ws.Range("E23").Value = '  -4.17%  '

# Row 24
$ws.Range("E24").Value = '  -0.27%  '

# Row 25
$ws.Range("D25").Value = '''150.24'
$ws.Range("E25").Value = '  -3.95%  '

# Row 26
$ws.Range("E26").Value = '  +1.94%  '

# Row 27
$ws.Range("E27").Value = '  -4.72%  '

# Row 28
$ws.Range("D28").Value = '''16.14'
$ws.Range("E28").Value = '  -2.37%  '

# Row 29
$ws.Range("D29").Value = '''1.383'
$ws.Range("E29").Value = '  -2.54%  '

# Row 30
$ws.Range("D30").Value = '''0.06131'
$ws.Range("E30").Value = '  -5.00%  '

# Row 31
$ws.Range("E31").Value = '  -3.72%  '

# Row 32
$ws.Range("D32").Value = '''3.736'
$ws.Range("E32").Value = '  -2.59%  '

# Row 33
$ws.Range("D33").Value = '''3.725'
$ws.Range("E33").Value = '  -1.93%  '

# Row 34
$ws.Range("D34").Value = '''1.670'
$ws.Range("E34").Value = '  -3.44%  '

# Row 35
$ws.Range("E35").Value = '  -5.69%  '

# Row 36
$ws.Range("D36").Value = '''0.6410'
$ws.Range("E36").Value = '  -1.14%  '

# Row 37
$ws.Range("D37").Value = '''2.414'
$ws.Range("E37").Value = '  -5.21%  '

# Row 38
$ws.Range("D38").Value = '''2.718'
$ws.Range("E38").Value = '  -1.36%  '

# Row 39
$ws.Range("D39").Value = '''0.01667'
$ws.Range("E39").Value = '  -5.14%  '

# Row 40
$ws.Range("D40").Value = '1.127.36'
$ws.Range("E40").Value = '  -1.25%  '

# Row 41
$ws.Range("D41").Value = '''6.159'
$ws.Range("E41").Value = '  -6.82%  '

# Row 42
$ws.Range("D42").Value = '''0.8723'
$ws.Range("E42").Value = '  -2.60%  '

# Row 43
$ws.Range("D43").Value = '''1.003'
$ws.Range("E43").Value = '  -0.19%  '

# Row 44
$ws.Range("D44").Value = '''99.62'
$ws.Range("E44").Value = '  -1.14%  '

# Row 45
$ws.Range("D45").Value = '1.887.19'
$ws.Range("E45").Value = '  -5.36%  '

# Row 46
$ws.Range("E46").Value = '  -4.78%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '''0.00000000108'
$ws.Range("E47").Value = '  -4.65%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.578'
$ws.Range("E48").Value = '  -2.89%  '

# Row 49
$ws.Range("D49").Value = '''8.274'
$ws.Range("E49").Value = '  -2.55%  '

# Row 50
$ws.Range("D50").Value = '''0.05380'
$ws.Range("E50").Value = '  -2.18%  '

# Row 51
$ws.Range("D51").Value = '''0.4414'
$ws.Range("E51").Value = '  -2.80%  '
